$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-09-06 Saturday" "2025-09-07 Sunday"

Replace-Text "96÷8=" "91÷6="
Replace-Text "90÷7=" "99÷9="
Replace-Text "81÷8=" "45÷2="
Replace-Text "54÷4=" "76÷3="
Replace-Text "54÷7=" "13÷9="

Replace-Text "32÷9=" "99÷9="
Replace-Text "38÷8=" "80÷4="
Replace-Text "88÷6=" "21÷6="
Replace-Text "26÷6=" "70÷5="
Replace-Text "40÷5=" "77÷3="

Replace-Text "18÷9=" "49÷6="
Replace-Text "26÷2=" "34÷7="
Replace-Text "17÷8=" "15÷8="
Replace-Text "74÷5=" "79÷6="
Replace-Text "60÷4=" "32÷7="

Replace-Text "57÷7=" "88÷7="
Replace-Text "91÷9=" "37÷7="
Replace-Text "26÷9=" "38÷4="
Replace-Text "52÷9=" "37÷6="
Replace-Text "41÷5=" "45÷3="

Replace-Text "72÷5=" "63÷4="
Replace-Text "60÷3=" "39÷9="
Replace-Text "91÷8=" "47÷9="
Replace-Text "26÷8=" "35÷5="
Replace-Text "17÷7=" "53÷5="
